# Added filtering options for the Component Analysis
#
# The workbook contains a "forecast-error triangle": rows are forecast
# origin dates, columns B:K are forecast horizons Q0..Q9. Values beyond
# what can be realized (i.e. cells that fall past the latest available
# actuals) must be filtered out / blanked.  This edit clears that set of
# now-out-of-range cells for rows 2-39 so the workbook matches the newly
# filtered component analysis output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$clearAddress = "I2:K2,G3:K3,I4:K4,G5:K5,I6:K6,G7:K7,I8:K8,G9:K9,I10:K10,G11:K11,I12:K12,K13,I14:K14,K15,I16:K16,K17,J18:K18,I19:K19,K21,J22:K22,I23:K23,K25,J26:K26,I27:K27,K29,J30:K30,I31:K31,K33,J34:K34,I35:K35,K37,J38:K38,I39:J39"

# NOTE: calling ClearContents() directly on a multi-area Range only
# affects the first area in this runtime, so iterate each area instead.
$rng = $ws.Range($clearAddress)
foreach ($area in $rng.Areas) {
    $area.ClearContents()
}
